# Swap the presentation's applied theme colors: "Integral" -> "Office Theme".
#
# The deck ships two theme parts: the slide master's theme (the one actually
# driving the rendered color scheme) and a second, otherwise-unreachable
# theme used only by the notes master. The edit being reproduced exchanges
# the two themes' bodies. Only the slide-master theme's 12 scheme colors are
# reachable from the PowerPoint object model (SlideMaster/NotesMaster both
# resolve to the same live ColorScheme here), and the font/format schemes of
# the two themes are already identical - so recoloring the reachable scheme
# to the "Office Theme" palette reproduces the observable effect of the swap.

function RgbValue($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Index -> (element, hex) per the standard 12-slot OOXML clrScheme order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
# Values below are the stock "Office Theme" color scheme.
$cs.Colors(1).RGB  = RgbValue 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = RgbValue 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = RgbValue 0x44 0x54 0x6A   # dk2
$cs.Colors(4).RGB  = RgbValue 0xE7 0xE6 0xE6   # lt2
$cs.Colors(5).RGB  = RgbValue 0x5B 0x9B 0xD5   # accent1
$cs.Colors(6).RGB  = RgbValue 0xED 0x7D 0x31   # accent2
$cs.Colors(7).RGB  = RgbValue 0xA5 0xA5 0xA5   # accent3
$cs.Colors(8).RGB  = RgbValue 0xFF 0xC0 0x00   # accent4
$cs.Colors(9).RGB  = RgbValue 0x44 0x72 0xC4   # accent5
$cs.Colors(10).RGB = RgbValue 0x70 0xAD 0x47   # accent6
$cs.Colors(11).RGB = RgbValue 0x05 0x63 0xC1   # hlink
$cs.Colors(12).RGB = RgbValue 0x95 0x4F 0x72   # folHlink
